$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the date labels in column A (rows 2-20) from "YYYY-MM" to "YYYY/MM".
for ($r = 2; $r -le 20; $r++) {
    $cur = $ws.Cells.Item($r, 1).Text
    $new = $cur -replace '-', '/'
    $ws.Cells.Item($r, 1).Value = $new
}

# 2. Normalize row heights for rows 7-20 to match rows 1-6 (13.8).
for ($r = 7; $r -le 20; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# 3. Update the current selection to A14:A20 with active cell A14.
$ws.Range("A14:A20").Select() | Out-Null
